# Updated symbol list with GitHub Actions
# Refreshes Price (D), Volume(1h) (E) and Hora (G) for each changed coin row.
# The source sheet stores these as plain text (inline strings), so every
# write forces the cell back to Text format first -- otherwise Excel would
# silently reinterpret "247.23" / "1.20%" / "13" as a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = '247.23'; Volume = '1.20%'; Hora = '13' }
    @{ Row = 3; Price = '30.01'; Volume = '11.00%'; Hora = '13' }
    @{ Row = 4; Price = '5.165'; Volume = '0.14%'; Hora = '13' }
    @{ Row = 5; Price = '0.05723'; Volume = '1.21%'; Hora = '13' }
    @{ Row = 6; Price = '6.617'; Volume = '2.20%'; Hora = '13' }
    @{ Row = 7; Price = '3.049'; Volume = '1.42%'; Hora = '13' }
    @{ Row = 8; Price = '0.8607'; Volume = '5.02%'; Hora = '13' }
    @{ Row = 9; Price = '0.8681'; Volume = '2.82%'; Hora = '13' }
    @{ Row = 10; Price = '0.1365'; Volume = '2.76%'; Hora = '13' }
    @{ Row = 11; Price = '0.07107'; Volume = '2.73%'; Hora = '13' }
    @{ Row = 12; Price = '0.02865'; Volume = '-2.36%'; Hora = '13' }
    @{ Row = 13; Price = '0.09380'; Volume = '-0.23%'; Hora = '13' }
    @{ Row = 14; Price = '0.001522'; Volume = '-0.08%'; Hora = '13' }
    @{ Row = 15; Price = '0.04140'; Volume = '-3.47%'; Hora = '13' }
    @{ Row = 16; Price = '0.0005989'; Volume = '0.02%'; Hora = '13' }
    @{ Row = 17; Price = '0.006179'; Volume = '0.31%'; Hora = '13' }
    @{ Row = 18; Price = '3.482'; Volume = '-0.79%'; Hora = '13' }
    @{ Row = 19; Price = '2.180'; Volume = '-2.11%'; Hora = '13' }
    @{ Row = 20; Price = '0.3190'; Volume = '2.45%'; Hora = '13' }
    @{ Row = 21; Price = '0.03245'; Volume = '3.07%'; Hora = '13' }
    @{ Row = 22; Price = $null; Volume = '4.09%'; Hora = '13' }
    @{ Row = 23; Price = '3.144'; Volume = '-12.51%'; Hora = '13' }
    @{ Row = 24; Price = '0.1380'; Volume = '0.46%'; Hora = '13' }
    @{ Row = 25; Price = '0.005108'; Volume = '14.24%'; Hora = '13' }
    @{ Row = 26; Price = '0.001220'; Volume = '-0.09%'; Hora = '13' }
    @{ Row = 27; Price = $null; Volume = '23.49%'; Hora = '13' }
    @{ Row = 28; Price = '0.0001938'; Volume = '167.06%'; Hora = '13' }
    @{ Row = 29; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 30; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 31; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 32; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 33; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 34; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 35; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 36; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 37; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 38; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 39; Price = $null; Volume = $null; Hora = '13' }
    @{ Row = 40; Price = '0.03780'; Volume = '3.64%'; Hora = '13' }
    @{ Row = 41; Price = '0.005959'; Volume = '-1.63%'; Hora = '13' }
    @{ Row = 42; Price = '0.1073'; Volume = '1.95%'; Hora = '13' }
    @{ Row = 43; Price = '0.002599'; Volume = '44.47%'; Hora = '13' }
    @{ Row = 44; Price = '0.009800'; Volume = '17.10%'; Hora = '13' }
    @{ Row = 45; Price = '0.00005090'; Volume = '-5.24%'; Hora = '13' }
    @{ Row = 46; Price = $null; Volume = '0.03%'; Hora = '13' }
    @{ Row = 47; Price = '0.07499'; Volume = '-31.81%'; Hora = '13' }
    @{ Row = 48; Price = '0.002761'; Volume = '3.87%'; Hora = '13' }
    @{ Row = 49; Price = '0.00002100'; Volume = '0.03%'; Hora = '13' }
    @{ Row = 50; Price = '0.0002000'; Volume = '0.03%'; Hora = '13' }
    @{ Row = 51; Price = $null; Volume = $null; Hora = '13' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $cell = $ws.Cells.Item($u.Row, 4)   # column D - Price
        $cell.NumberFormat = "@"
        $cell.Value = $u.Price
    }
    if ($null -ne $u.Volume) {
        $cell = $ws.Cells.Item($u.Row, 5)   # column E - Volume(1h)
        $cell.NumberFormat = "@"
        $cell.Value = $u.Volume
    }
    if ($null -ne $u.Hora) {
        $cell = $ws.Cells.Item($u.Row, 7)   # column G - Hora
        $cell.NumberFormat = "@"
        $cell.Value = $u.Hora
    }
}

